# Generate Report for Handback
# Refresh the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the 34493f90-67fa-4c88-939f-d076ab8bf2b3 row, as produced by a re-run
# of the handback report generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (row for 34493f90-...)
$wsOverview.Range("G3").Value = "2016-08-20 08:53:29"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H3").Value = "2016-08-20 08:53:25"
$wsZhCn.Range("K3").Value = "2016-08-20 08:53:41"

# de-de sheet: Correspond Handoff Datetime (mirrors Overview's HO Xliff date)
# and Correspond Handback DateTime
$wsDeDe.Range("H3").Value = "2016-08-20 08:53:29"
$wsDeDe.Range("K3").Value = "2016-08-20 08:53:48"
